# Edit script: add 2022-Q4 data
# 1) Update "总计" (summary) sheet: insert a new row for 2022-Q4 at the top of
#    the data (row 2), pushing the existing quarters down by one row.
# 2) Insert a brand-new worksheet named "2022-Q4" right before the existing
#    "2022-Q3" sheet, populated with the fund-holding breakdown for that
#    quarter.
# All other quarter sheets (2022-Q3 .. 2021-Q1) are left untouched.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) 总计 (summary) sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Data rows currently occupy r2:D8 (quarter, count, value). We need to shift
# them down by one row (to r3:D9) and then write the new 2022-Q4 row at r2.
# Walk from the bottom up so we don't clobber rows before they are read.
for ($r = 8; $r -ge 2; $r--) {
    $srcA = $summary.Cells.Item($r, 1).Value2
    $srcB = $summary.Cells.Item($r, 2).Value2
    $srcC = $summary.Cells.Item($r, 3).Value2
    $srcD = $summary.Cells.Item($r, 4).Value2

    $dst = $r + 1
    $summary.Cells.Item($dst, 1).Value = $srcA + 1
    $summary.Cells.Item($dst, 2).NumberFormat = "@"
    $summary.Cells.Item($dst, 2).Value = $srcB
    $summary.Cells.Item($dst, 3).Value = $srcC
    $summary.Cells.Item($dst, 4).Value = $srcD
}

# New top row: 2022-Q4
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).NumberFormat = "@"
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 11
$summary.Cells.Item(2, 4).Value = 3.61

# ---------------------------------------------------------------------
# 2) New "2022-Q4" sheet (inserted before "2022-Q3")
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($anchor)
$q4.Name = "2022-Q4"

# Header row
$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

# Helper to write one data row: index, code, name, size, position, weight, value, rank
function Write-FundRow($sheet, $row, $idx, $code, $name, $size, $pos, $weight, $value, $rank) {
    $sheet.Cells.Item($row, 1).Value = $idx

    $sheet.Cells.Item($row, 2).NumberFormat = "@"
    $sheet.Cells.Item($row, 2).Value = $code

    $sheet.Cells.Item($row, 3).NumberFormat = "@"
    $sheet.Cells.Item($row, 3).Value = $name

    $sheet.Cells.Item($row, 4).NumberFormat = "@"
    $sheet.Cells.Item($row, 4).Value = $size

    $sheet.Cells.Item($row, 5).NumberFormat = "@"
    $sheet.Cells.Item($row, 5).Value = $pos

    $sheet.Cells.Item($row, 6).NumberFormat = "@"
    $sheet.Cells.Item($row, 6).Value = $weight

    $sheet.Cells.Item($row, 7).NumberFormat = "@"
    $sheet.Cells.Item($row, 7).Value = $value

    $sheet.Cells.Item($row, 8).Value = $rank
}

Write-FundRow $q4 2  0  "010387" "易方达医药生物股票A"                  "20.77" "91.88" "4.84" "1.0053" 6
Write-FundRow $q4 3  1  "012346" "易方达港股通成长混合A"                "30.88" "89.23" "3.03" "0.9357" 9
Write-FundRow $q4 4  2  "012347" "易方达港股通成长混合C"                "21.55" "89.23" "3.03" "0.6530" 9
Write-FundRow $q4 5  3  "007718" "中银创新医疗混合A"                    "11.61" "80.72" "3.43" "0.3982" 9
Write-FundRow $q4 6  4  "010388" "易方达医药生物股票C"                  "6.05"  "91.88" "4.84" "0.2928" 6
Write-FundRow $q4 7  5  "470888" "汇添富香港优势精选混合（QDII）"        "2.68"  "93.08" "4.92" "0.1319" 5
Write-FundRow $q4 8  6  "010500" "中银创新医疗混合C"                    "3.82"  "80.72" "3.43" "0.1310" 9
Write-FundRow $q4 9  7  "012086" "博时健康生活混合A"                    "2.45"  "93.31" "1.44" "0.0353" 10
Write-FundRow $q4 10 8  "008861" "西部利得港股通新机遇灵活配置混合A"     "0.25"  "87.69" "3.66" "0.0092" 4
Write-FundRow $q4 11 9  "012087" "博时健康生活混合C"                    "0.57"  "93.31" "1.44" "0.0082" 10
Write-FundRow $q4 12 10 "010093" "西部利得港股通新机遇灵活配置混合C"     "0.12"  "87.69" "3.66" "0.0044" 4
